$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data fix: two motion filenames ("0071@000_0_0.npy" and "0071@002_50_0.npy")
# were merged into a single renamed entry, shifting every following filename
# up by one cell. ---
$ws.Range("A2").Value = "s_07_act_08_subact_01_ca_01.npy"
$ws.Range("B2").Value = "s_07_act_08_subact_01_ca_01.npy"
$ws.Range("C2").Value = "0071@012_60_0.npy"
$ws.Range("D2").Value = "0071@019_0_0.npy"
$ws.Range("A3").Value = "0071@022_40_0.npy"
$ws.Range("B3").Value = "0071@027_60_1.npy"
$ws.Range("C3").Value = "0071@028_20_0.npy"
$ws.Range("D3").Value = "0071@035_10_1.npy"
$ws.Range("A4").Value = "0071@038_10_0.npy"
$ws.Range("B4").Value = "0071@039_10_1.npy"
$ws.Range("C4").Value = "0071@044_0_0.npy"
$ws.Range("D4").Value = "0071@045_60_0.npy"
$ws.Range("A5").Value = "0071@046_20_0.npy"
$ws.Range("B5").Value = "0071@050_10_0.npy"
$ws.Range("C5").Value = "0071@056_30_1.npy"
$ws.Range("D5").Value = "0071@061_30_1.npy"

# --- Row 1 grew slightly taller ---
$ws.Rows.Item(1).RowHeight = 20.25

# --- C1/D1 (currently blank, unformatted cells next to the "bpm" header) now
# pick up the same border + a plain Calibri font, consistent with the rest
# of the header row. ---
$ws.Range("B1").Copy()
$ws.Range("C1:D1").PasteSpecial(-4122)
$ws.Range("C1:D1").Font.ThemeColor = 1
